$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns; rows 43-51
# additionally had their Coin (B) / Link (C) values re-ranked.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.793.81'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.37%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.633.77'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.78%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.55'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.50%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.569'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.60'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -6.37%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.103'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.335'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.099.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '58.781.01'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.34%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.92'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000136'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.652.32'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '340.24'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.45'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.49'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.32'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.418'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.21'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0795'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.15%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.46'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.74%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.89'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '149.45'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.19'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.19'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.866'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.80%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '36.41'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('E39').Value = '  -4.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.64'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.604'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0974'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '270.41'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.39%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.65'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0537'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.26'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.037.81'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.72%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.75'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.44%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0230'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.94'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.42%  '
